# Updates cryptos list (prices / 1h volume %) to the latest scraped values.
# Leading apostrophes on some D-column (Price) values force Excel to keep
# storing them as text (matching the source data, which is plain text),
# instead of auto-converting number-looking strings like "351.77" into
# numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.798.19'
$ws.Range('E2').Value = '  +1.90%  '
$ws.Range('D3').Value = '2.806.00'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''351.77'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').Value = '''112.29'
$ws.Range('E6').Value = '  +5.33%  '
$ws.Range('D7').Value = '''0.557'
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '''0.621'
$ws.Range('E9').Value = '  +7.68%  '
$ws.Range('D10').Value = '''40.20'
$ws.Range('E10').Value = '  +3.27%  '
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').Value = '''0.0838'
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('D13').Value = '''19.91'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').Value = '''7.77'
$ws.Range('E14').Value = '  +4.43%  '
$ws.Range('D15').Value = '3.244.00'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').Value = '2.802.47'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '''0.959'
$ws.Range('E17').Value = '  +4.00%  '
$ws.Range('D18').Value = '51.798.76'
$ws.Range('E18').Value = '  +2.00%  '
$ws.Range('B19').Value = 'ImmutableX'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D19').Value = '''3.32'
$ws.Range('E19').Value = '  +9.71%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '''7.62'
$ws.Range('E20').Value = '  +1.57%  '
$ws.Range('D21').Value = '''13.53'
$ws.Range('E21').Value = '  +4.64%  '
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('D23').Value = '''70.35'
$ws.Range('E23').Value = '  +1.73%  '
$ws.Range('D24').Value = '''267.62'
$ws.Range('E24').Value = '  +1.73%  '
$ws.Range('D25').Value = '''2.75'
$ws.Range('E25').Value = '  +2.33%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = '''26.18'
$ws.Range('E27').Value = '  +1.76%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '''39.11'
$ws.Range('E29').Value = '  +14.50%  '
$ws.Range('E30').Value = '  +3.85%  '
$ws.Range('E31').Value = '  +2.40%  '
$ws.Range('D32').Value = '''52.58'
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('D33').Value = '''6.11'
$ws.Range('E33').Value = '  +1.81%  '
$ws.Range('D34').Value = '''0.0902'
$ws.Range('E34').Value = '  +9.36%  '
$ws.Range('D35').Value = '''0.0451'
$ws.Range('E35').Value = '  +1.98%  '
$ws.Range('D36').Value = '''5.56'
$ws.Range('E36').Value = '  +4.76%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').Value = '''19.01'
$ws.Range('E38').Value = '  +4.58%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '''3.17'
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '''2.01'
$ws.Range('E40').Value = '  +3.91%  '
$ws.Range('E41').Value = '  +2.56%  '
$ws.Range('D42').Value = '''2.51'
$ws.Range('E42').Value = '  +2.82%  '
$ws.Range('E43').Value = '  +1.92%  '
$ws.Range('D44').Value = '''120.59'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').Value = '''21.98'
$ws.Range('E45').Value = '  +0.46%  '
$ws.Range('D46').Value = '''3.53'
$ws.Range('E46').Value = '  +10.19%  '
$ws.Range('D47').Value = '''2.47'
$ws.Range('E47').Value = '  +8.81%  '
$ws.Range('D48').Value = '2.122.57'
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('D49').Value = '''0.979'
$ws.Range('E49').Value = '  +7.94%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '''5.49'
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '''1.37'
$ws.Range('E51').Value = '  +7.92%  '
